# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.570.72'
$ws.Range("E2").Value = '  -0.46%  '

$ws.Range("D3").Value = '3.169.13'
$ws.Range("E3").Value = '  -5.84%  '

$ws.Range("E4").Value = '  -0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.51'
$ws.Range("E5").Value = '  -4.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '608.63'
$ws.Range("E6").Value = '  -5.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.377'
$ws.Range("E7").Value = '  -8.52%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.664'
$ws.Range("E8").Value = '  +0.17%  '

$ws.Range("E9").Value = '  -0.17%  '

$ws.Range("D10").Value = '3.163.22'
$ws.Range("E10").Value = '  -5.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.534'
$ws.Range("E11").Value = '  -12.85%  '

$ws.Range("E12").Value = '  +6.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000241'
$ws.Range("E13").Value = '  -14.08%  '

$ws.Range("D14").Value = '3.750.64'
$ws.Range("E14").Value = '  -6.34%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.23'
$ws.Range("E15").Value = '  -4.47%  '

$ws.Range("D16").Value = '87.517.93'
$ws.Range("E16").Value = '  -0.54%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '32.11'
$ws.Range("E17").Value = '  -10.73%  '

$ws.Range("D18").Value = '3.185.77'
$ws.Range("E18").Value = '  -5.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.08'
$ws.Range("E19").Value = '  +2.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.38'
$ws.Range("E20").Value = '  -9.66%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '413.29'
$ws.Range("E21").Value = '  -8.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.43'
$ws.Range("E22").Value = '  -11.69%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.04'
$ws.Range("E23").Value = '  -8.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.16'
$ws.Range("E24").Value = '  -4.88%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.82'
$ws.Range("E25").Value = '  -4.38%  '

$ws.Range("D26").Value = '3.327.63'
$ws.Range("E26").Value = '  -6.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '73.02'
$ws.Range("E27").Value = '  -8.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000130'
$ws.Range("E28").Value = '  -7.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.05%  '

$ws.Range("E30").Value = '  -0.08%  '

$ws.Range("E31").Value = '  -14.92%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '544.15'
$ws.Range("E32").Value = '  -4.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.19'
$ws.Range("E33").Value = '  -11.35%  '

$ws.Range("E34").Value = '  -15.16%  '

$ws.Range("E35").Value = '  -11.01%  '

$ws.Range("E36").Value = '  -8.24%  '

$ws.Range("E37").Value = '  -6.36%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '21.76'
$ws.Range("E38").Value = '  -7.66%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '21.83'
$ws.Range("E39").Value = '  -0.15%  '

$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.00'
$ws.Range("E41").Value = '  -2.52%  '

$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("E43").Value = '  -8.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.367'
$ws.Range("E44").Value = '  -14.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.82'
$ws.Range("E45").Value = '  -5.48%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '172.27'
$ws.Range("E46").Value = '  -7.50%  '

$ws.Range("E47").Value = '  -7.18%  '

$ws.Range("E48").Value = '  +1.44%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.22'
$ws.Range("E49").Value = '  -13.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.95'
$ws.Range("E50").Value = '  -11.29%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.55'
$ws.Range("E51").Value = '  -6.81%  '
